# Update "Pais" (countries) COVID data sheet:
#  - refresh the "last updated" timestamp
#  - update stat figures for a set of countries
#  - Indonesia overtakes Ecuador in total cases -> rows swap (ranking by column B)
#  - Fiyi/Dominica and Islas Malvinas/Groenlandia swap their tied ranking order

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 12:15"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2682011
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 1122246
$ws.Range("E4").Value = 1430977
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 128788

# Row 7 - India
$ws.Range("B7").Value = 568346
$ws.Range("C7").Value = 810
$ws.Range("D7").Value = 335656
$ws.Range("E7").Value = 215771
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 16919

# Row 17 - Alemania
$ws.Range("B17").Value = 195399
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 7258

# Row 30 - Belgica
$ws.Range("B30").Value = 61427
$ws.Range("C30").Value = 66
$ws.Range("D30").Value = 16984
$ws.Range("E30").Value = 34696
$ws.Range("G30").Value = 15
$ws.Range("H30").Value = 9747

# Row 31 - was Ecuador, now Indonesia overtakes it (new, higher total cases)
$ws.Range("A31").Value = "Indonesia"
$ws.Range("B31").Value = 56385
$ws.Range("C31").Value = 1293
$ws.Range("D31").Value = 24806
$ws.Range("E31").Value = 28703
$ws.Range("G31").Value = 71
$ws.Range("H31").Value = 2876

# Row 32 - was Indonesia, now Ecuador (its figures are unchanged, just shifted down a rank)
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 55665
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 27430
$ws.Range("E32").Value = 23733
$ws.Range("H32").Value = 4502

# Row 40 - Oman
$ws.Range("B40").Value = 40070
$ws.Range("C40").Value = 1010
$ws.Range("D40").Value = 23425
$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 176

# Row 41 - Filipinas
$ws.Range("B41").Value = 37514
$ws.Range("C41").Value = 1076
$ws.Range("D41").Value = 10233
$ws.Range("E41").Value = 26015
$ws.Range("G41").Value = 11
$ws.Range("H41").Value = 1266

# Row 49 - Barein
$ws.Range("E49").Value = 5225
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 86

# Row 61 - Moldavia
$ws.Range("D61").Value = 9382
$ws.Range("E61").Value = 6434
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 541

# Row 68 - Marruecos
$ws.Range("B68").Value = 12385
$ws.Range("C68").Value = 95
$ws.Range("D68").Value = 8839
$ws.Range("E68").Value = 3321

# Row 73 - Malasia
$ws.Range("B73").Value = 8639
$ws.Range("C73").Value = 2
$ws.Range("D73").Value = 8354
$ws.Range("E73").Value = 164

# Row 103 - Albania
$ws.Range("B103").Value = 2535
$ws.Range("C103").Value = 69
$ws.Range("D103").Value = 1459
$ws.Range("E103").Value = 1014
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 62

# Rows 205/206 - Dominica and Fiyi are tied; swap their displayed order
$ws.Range("A205").Value = "Fiyi"
$ws.Range("A206").Value = "Dominica"

# Rows 209/210 - Groenlandia and Islas Malvinas are tied; swap their displayed order
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
